# Convert the NCSTE publication-certificate letter template so the
# {#publications}...{/publications} Mustache/Carbone loop wraps the whole
# table (as a table-per-publication block) instead of being repeated
# individually inside every cell.

$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# 1. Un-wrap each table cell's per-field loop: {#publications}{field}{/publications} -> {field}
$d.Content.Find.Execute("{#publications}{no}. Название статьи:{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{no}. Название статьи:", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{title}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{title}", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{authors}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{authors}", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{journal}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{journal}", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{volume_issue}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{volume_issue}", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{vol}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{vol}", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{year}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{year}", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{issn_print}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{issn_print}", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{issn_online}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{issn_online}", 2) | Out-Null

$d.Content.Find.Execute("{#publications}{doi}{/publications}", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "{doi}", 2) | Out-Null

# 2. Insert a new paragraph containing the opening "{#publications}" tag
#    right before the table (the empty paragraph with spacing-after=100,
#    which directly precedes the table, is paragraph #14).
$beforeTablePara = $d.Paragraphs.Item(14)
$beforeTableRange = $beforeTablePara.Range
$openTagXml = '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:spacing w:after="100"/></w:pPr></w:p>' + `
              '<w:p xmlns:w="' + $wNs + '"><w:r><w:t>{#publications}</w:t></w:r></w:p>'
$beforeTableRange.InsertXML($openTagXml)

# 3. Insert a new paragraph containing the closing "{/publications}" tag
#    right after the table. After step 2 the empty paragraph (spacing-after
#    =200) that follows the table is now paragraph #43.
$afterTablePara = $d.Paragraphs.Item(43)
$afterTableRange = $afterTablePara.Range
$closeTagXml = '<w:p xmlns:w="' + $wNs + '"><w:r><w:t>{/publications}</w:t></w:r></w:p>' + `
               '<w:p xmlns:w="' + $wNs + '"><w:pPr><w:spacing w:after="200"/></w:pPr></w:p>'
$afterTableRange.InsertXML($closeTagXml)
